# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) for the first handback entry
# (23efbcc1-...) on both the zh-cn and de-de status sheets, reflecting a
# new handback report run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-11 18:32:46"
$zhcn.Range("H2").Value = "2016-03-11 18:33:07"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-11 18:32:49"
$dede.Range("H2").Value = "2016-03-11 18:33:12"
